# "Another round of runs" — refresh the benchmark numbers on Sheet1 with
# the latest set of measurements (service-compare-data/Data.xlsx).
#
# The workbook's chart sheets (GetChart / SendChart / FirstChart /
# SendOneChart) read their cached series data from these cells; Excel will
# refresh the chart caches itself the next time it opens/recalculates the
# workbook, so only the underlying Sheet1 values need to be updated here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Web API (Core 3.1)
$ws.Range("B6").Value = 33.9
$ws.Range("C6").Value = 3.16
$ws.Range("D6").Value = 80.2
$ws.Range("E6").Value = 15.7
$ws.Range("F6").Value = 32.4
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 52.1
$ws.Range("I6").Value = 6.24
$ws.Range("J6").Value = 0.99
$ws.Range("K6").Value = 0.28000000000000003
$ws.Range("L6").Value = 1.98
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 7.1
$ws.Range("O6").Value = 0.83
$ws.Range("P6").Value = 15.5
$ws.Range("Q6").Value = 4.6399999999999997

# gRPC (Core 3.1)
$ws.Range("B7").Value = 21.1
$ws.Range("C7").Value = 1.77
$ws.Range("D7").Value = 66.2
$ws.Range("E7").Value = 19.100000000000001
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = 2.33
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 16.399999999999999
$ws.Range("J7").Value = 1.07
$ws.Range("K7").Value = 0.18
$ws.Range("L7").Value = 2.37
$ws.Range("M7").Value = 2.27
$ws.Range("N7").Value = 0.94
$ws.Range("O7").Value = 0.18
$ws.Range("P7").Value = 1.84
$ws.Range("Q7").Value = 1.94

# gRPC Streaming (Core 3.1)
$ws.Range("B8").Value = 73.400000000000006
$ws.Range("C8").Value = 5.12
$ws.Range("D8").Value = 435
$ws.Range("E8").Value = 228

# Web Api (Framework 4.8)
$ws.Range("B9").Value = 54.6
$ws.Range("C9").Value = 3.36
$ws.Range("D9").Value = 104
$ws.Range("E9").Value = 22.6
$ws.Range("F9").Value = 101
$ws.Range("G9").Value = 16.899999999999999
$ws.Range("H9").Value = 240
$ws.Range("I9").Value = 30
$ws.Range("J9").Value = 1.52
$ws.Range("K9").Value = 0.42
$ws.Range("L9").Value = 2.39
$ws.Range("M9").Value = 2.4700000000000002
$ws.Range("N9").Value = 7.91
$ws.Range("O9").Value = 1.1200000000000001
$ws.Range("P9").Value = 15.7
$ws.Range("Q9").Value = 5.09

# WCF (Framework 4.8)
$ws.Range("B10").Value = 72.3
$ws.Range("C10").Value = 16.399999999999999
$ws.Range("D10").Value = 146
$ws.Range("E10").Value = 31
$ws.Range("F10").Value = 75.099999999999994
$ws.Range("G10").Value = 12.2
$ws.Range("H10").Value = 160
$ws.Range("I10").Value = 23.2
$ws.Range("J10").Value = 3.06
$ws.Range("K10").Value = 2.2000000000000002
$ws.Range("L10").Value = 3.3
$ws.Range("M10").Value = 3.03
$ws.Range("N10").Value = 4.82
$ws.Range("O10").Value = 5
$ws.Range("P10").Value = 3.1
$ws.Range("Q10").Value = 2.56

# The author's last click before saving landed on A10 (selection moved down
# one row from the previous A9) and that sheet tab became the active one.
$ws.Range("A10").Select()
